# "Generate Report for Handback" - refresh the localization-status report
# after a successful handback run for zh-cn and de-de.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Status column updates (was "Ready for handoff") ---
$overview.Range("E2").Value = $statusText   # zh-cn status on the Overview sheet
$overview.Range("F2").Value = $statusText   # de-de status on the Overview sheet
$zhcn.Range("C2").Value = $statusText
$dede.Range("C2").Value = $statusText

# --- Latest Handback DateTime refreshed by this handback run ---
$zhcn.Range("K2").Value = "2016-08-16 16:45:21"
$dede.Range("K2").Value = "2016-08-16 16:45:28"

# --- Error Detail cleared now that handback is in sync with en-US ---
# (use the text-prefix trick so the cell keeps an explicit, empty
#  string value instead of being deleted outright)
$zhcn.Range("P2").Value = "'"
$zhcn.Range("P2").Style = "Normal"
$dede.Range("P2").Value = "'"
$dede.Range("P2").Style = "Normal"

# --- Column widths widened/narrowed to fit the refreshed report text ---
$overview.Columns.Item(5).ColumnWidth = 29.16666667
$overview.Columns.Item(6).ColumnWidth = 29.16666667
$zhcn.Columns.Item(3).ColumnWidth = 29.16666667
$dede.Columns.Item(3).ColumnWidth = 29.16666667
$zhcn.Columns.Item(16).ColumnWidth = 12.83333333
$dede.Columns.Item(16).ColumnWidth = 12.83333333
